$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.300.05'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '2.631.38'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.30'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.48'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.96'
$ws.Range('E9').Value = '  +8.07%  '
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '3.101.21'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '59.242.34'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.15'
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.664.39'
$ws.Range('E16').Value = '  +1.28%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.47'
$ws.Range('E18').Value = '  +2.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.16'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.23'
$ws.Range('E20').Value = '  +0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.23'
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.33'
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.416'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.164'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.27'
$ws.Range('E27').Value = '  +0.29%  '
$ws.Range('D28').Value = '0.0₃0749'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.65'
$ws.Range('E30').Value = '  -2.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.88'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.80'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '150.96'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.97'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.13'
$ws.Range('E35').Value = '  +0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.841'
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.832'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.44'
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.60'
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '285.20'
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.598'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.74'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0538'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.08'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0940'
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0226'
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('D48').Value = '1.957.24'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.54'
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.22'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.52'
$ws.Range('E51').Value = '  +0.42%  '
